# DPLKKEU009-001 -> rename sheet/test-case + refresh scenario data for a
# new run (Kantor ID 216/PASURUAN, INV.DEP.PEN.012, BC001-O-22-08-00012).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet / TC_ID rename -------------------------------------------------
$ws.Name = "DPLKKEU009-001"
$ws.Range("B2").Value = "DPLKKEU009-001"

# --- Scenario description (D2): settlement reference number changes ------
$ws.Range("D2").Value = "Approval Settlement Transaksi Penempatan Deposito BANK BNI - JPU Jakarta Pusat DLK/3/3162 PRI01-Deposito Pasar Uang DOC -- Keuangan Investasi - Proses - Approval Pembayaran Investasi"

# --- Preparation block (F2): new search/kode pembayaran/kantor/date/amount
$prep = "Username : 31816;`n" +
        "Password : bni1234;`n" +
        "Role : 09;`n" +
        "Search : INV.DEP.PEN.012;`n" +
        "Kode Pembayaran : BC001-O-22-08-00012;`n" +
        "Entitas : DPLKBNI : DPLK PT. BNI (Persero) Tbk.;`n" +
        "Kantor ID : 216 : PASURUAN;`n" +
        "Kode Buku : BC001;`n" +
        "Nama Bank : Bank Negara Indonesia 1946;`n" +
        "No Rekening : 1000564390;`n" +
        "Produk/Cluster : Deposito Pasar Uang;`n" +
        "Mata Uang : IDR : Rupiah;`n" +
        "Tanggal RK : 05/08/2022;`n" +
        "Nominal Pembayaran : 40.000.000.000,00;`n" +
        "Keterangan : INV.DEP.PEN.012 DPLK PT. BNI (Persero) Tbk. Deposito Pasar Uang"
$ws.Range("F2").Value = $prep

# --- KODE_PEMBAYARAN / KETERANGAN values ----------------------------------
$ws.Range("N2").Value = "BC001-O-22-08-00012"
$ws.Range("O2").Value = "INV.DEP.PEN.012"

# --- Stray empty/styled cell dropped by the editor on resave --------------
$ws.Range("M3").Clear()

# --- Selection / active cell as left by the editor ------------------------
$ws.Range("I1").Select()
